# "Satz > Besonderheit" - add a new "Tabelle2" sheet with category/relation
# lookup data + generated SQL INSERT statements, right after "Tabelle1".

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- create the new sheet right after Tabelle1 ---------------------------
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Tabelle2"

# --- column A (Name) : enter top-to-bottom first so shared-string ids ----
# --- line up with the authoring order captured in the workbook ----------
$ws2.Range("A2").Value  = "Verwendungszweck"
$ws2.Range("A3").Value  = "Besetzung"
$ws2.Range("A4").Value  = "Epoche"
$ws2.Range("A5").Value  = "Gattung"
$ws2.Range("A6").Value  = "Erprobt"
$ws2.Range("A7").Value  = "Notenwert"
$ws2.Range("A8").Value  = "Schwierigkeitsgrad"
$ws2.Range("A9").Value  = "Strichart"
$ws2.Range("A10").Value = "Übung"
$ws2.Range("A11").Value = "Melodische Besonderheit"
$ws2.Range("A12").Value = "Dynamische Besonderheit"
$ws2.Range("A13").Value = "Rhythmische Besonderheit "

# --- header row -----------------------------------------------------------
$ws2.Range("A1").Value = "Name"
$ws2.Range("B1").Value = "Relation"
$ws2.Range("C1").Value = "SQL "

# --- column B (Relation) ---------------------------------------------------
$ws2.Range("B2:B5").Value  = "musikstueck"
$ws2.Range("B6:B13").Value = "satz"

# --- column C (generated SQL) ---------------------------------------------
$ws2.Range("C2").Formula = '="INSERT INTO category (Name, Relation) VALUES(''" & A2 & "'', ''" & B2 & "'') ;"'
$ws2.Range("C3:C13").Formula = '="INSERT INTO category (Name, Relation) VALUES(''" & A3 & "'', ''" & B3 & "'') ;"'

# --- header formatting ------------------------------------------------------
$ws2.Range("A1:D1").Font.Bold = $true

# --- column widths -----------------------------------------------------------
$ws2.Columns.Item(1).ColumnWidth = 21.42
$ws2.Columns.Item(2).ColumnWidth = 21.42

# --- page setup (match Tabelle1) --------------------------------------------
$ps = $ws2.PageSetup
$ps.TopMargin = 56.6929134
$ps.BottomMargin = 56.6929134
$ps.LeftMargin = 50.4
$ps.RightMargin = 50.4
$ps.HeaderMargin = 21.6
$ps.FooterMargin = 21.6
$ps.PaperSize = 9
$ps.Orientation = 1

# --- selection / active sheet ----------------------------------------------
$ws1.Range("H48").Select() | Out-Null
$ws2.Range("E17").Select() | Out-Null
$ws2.Activate() | Out-Null
